$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.878.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.17%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.783.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.88%  '

# Row 4
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.43%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.46%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.780.62'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.98%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.448'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.71%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.82'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.11%  '

# Row 13
$ws.Range("E13").Value = '  -2.43%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.00%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.416.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.89%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.787.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.59%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.856.89'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.10%  '

# Row 18
$ws.Range("E18").Value = '  -1.23%  '

# Row 19
$ws.Range("E19").Value = '  +1.83%  '

# Row 20
$ws.Range("E20").Value = '  -1.01%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '459.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.84%  '

# Row 22
$ws.Range("E22").Value = '  -4.42%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.692'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.04%  '

# Row 24
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000146'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.58%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.13%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.01%  '

# Row 28
$ws.Range("E28").Value = '  -0.03%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.91'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.33%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.922.80'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.05%  '

# Row 31
$ws.Range("E31").Value = '  -6.72%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.07%  '

# Row 33
$ws.Range("E33").Value = '  -1.88%  '

# Row 34
$ws.Range("E34").Value = '  -1.96%  '

# Row 35
$ws.Range("E35").Value = '  -0.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.69%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0991'
$ws.Range("D37").Style = "Normal"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.148'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.48%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.81'
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.25%  '

# Row 41
$ws.Range("E41").Value = '  -2.09%  '

# Row 42
$ws.Range("E42").Value = '  +0.07%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.77%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.03%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '152.41'
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = '  -1.91%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.65%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.31'
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = '  -0.25%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.37'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.17%  '
